# scrum_v02.xlsx update
# - window height (cosmetic, not controllable through this COM surface -> left as-is)
# - sharedStrings grow with new Sprint Backlog content ("done" status + 8 new sprint items)
# - Product Backlog: selection + column C width
# - Sprint Backlog: selection, column D width, status of items 1.1-1.3 -> done,
#   and 8 new backlog rows (2.1 - 4.1 and 1.4) added

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Product Backlog sheet
# ---------------------------------------------------------------------------
$productBacklog = $wb.Worksheets.Item("Product Backlog")

# widen column C (closest achievable width to the target 38.44140625 through this
# COM surface, which always rounds column widths to whole pixels)
$productBacklog.Columns.Item(3).ColumnWidth = 37.666666666666664

$productBacklog.Range("B3:B5").Select()

# ---------------------------------------------------------------------------
# Sprint Backlog sheet
# ---------------------------------------------------------------------------
$sprintBacklog = $wb.Worksheets.Item("Sprint Backlog")
$sprintBacklog.Activate()

# widen column D (closest achievable width to the target 85.44140625)
$sprintBacklog.Columns.Item(4).ColumnWidth = 84.66666666666667

# mark the 3 existing sprint 1 items as done (remaining effort 12 -> 12, status -> done)
$sprintBacklog.Range("K2").Value = 12
$sprintBacklog.Range("L2").Value = "done"
$sprintBacklog.Range("K3").Value = 12
$sprintBacklog.Range("L3").Value = "done"
$sprintBacklog.Range("K4").Value = 12
$sprintBacklog.Range("L4").Value = "done"

# Row 5 - Medication entry View
$sprintBacklog.Range("A5").Value = 2.1
$sprintBacklog.Range("B5").Value = 2
$sprintBacklog.Range("C5").Value = "Medication entry View"
$sprintBacklog.Range("D5").Value = "Create the planned Views for the entry of a new Medicament"
$sprintBacklog.Range("E5").Value = "UI, Controller"
$sprintBacklog.Range("F5").Value = "Mete"
$sprintBacklog.Range("G5").Value = "Amin"
$sprintBacklog.Range("I5").Value = 20
$sprintBacklog.Range("L5").Value = "work in progress"
$sprintBacklog.Rows.Item(5).RowHeight = 14.4

# Row 6 - Medication Data Model
$sprintBacklog.Range("A6").Value = 2.2
$sprintBacklog.Range("B6").Value = 2
$sprintBacklog.Range("C6").Value = "Medication Data Model"
$sprintBacklog.Range("D6").Value = "Optimize the medication dose schema to enter a new Medicament"
$sprintBacklog.Range("E6").Value = "Modell, Database"
$sprintBacklog.Range("F6").Value = "Marwin"
$sprintBacklog.Range("G6").Value = "Michel"
$sprintBacklog.Range("I6").Value = 10
$sprintBacklog.Range("L6").Value = "work in progress"

# Row 7 - Medication validation
$sprintBacklog.Range("A7").Value = 2.3
$sprintBacklog.Range("B7").Value = 2
$sprintBacklog.Range("C7").Value = "Medication validation"
$sprintBacklog.Range("D7").Value = "Validate and save an entered medicament"
$sprintBacklog.Range("E7").Value = "Controller, Modell"
$sprintBacklog.Range("F7").Value = "Jonas"
$sprintBacklog.Range("G7").Value = "Marwin"
$sprintBacklog.Range("I7").Value = 15
$sprintBacklog.Range("L7").Value = "work in progress"

# Row 8 - Navigation
$sprintBacklog.Range("A8").Value = 2.4
$sprintBacklog.Range("B8").Value = 2
$sprintBacklog.Range("C8").Value = "Navigation"
$sprintBacklog.Range("D8").Value = "Navigate to the new View and Back"
$sprintBacklog.Range("E8").Value = "UI, Controller"
$sprintBacklog.Range("F8").Value = "Carole"
$sprintBacklog.Range("G8").Value = "Jonas"
$sprintBacklog.Range("I8").Value = 4
$sprintBacklog.Range("L8").Value = "work in progress"

# Row 9 - Medication edit View (wrapped, taller row)
$sprintBacklog.Range("A9").Value = 3.1
$sprintBacklog.Range("B9").Value = 2
$sprintBacklog.Range("C9").Value = "Medication edit View"
$sprintBacklog.Range("D9").Value = "Create the planned Editing view to edit an existing medicament by using 2,3.`nDo not implement history yet."
$sprintBacklog.Range("D9").WrapText = $true
$sprintBacklog.Range("E9").Value = "UI, Controller"
$sprintBacklog.Range("F9").Value = "Amin"
$sprintBacklog.Range("G9").Value = "Mete"
$sprintBacklog.Range("I9").Value = 10
$sprintBacklog.Range("L9").Value = "work in progress"
$sprintBacklog.Rows.Item(9).RowHeight = 33

# Row 10 - Navigation
$sprintBacklog.Range("A10").Value = 3.2
$sprintBacklog.Range("B10").Value = 2
$sprintBacklog.Range("C10").Value = "Navigation"
$sprintBacklog.Range("D10").Value = "Navigate to the new view over the Medicationoverview-List"
$sprintBacklog.Range("E10").Value = "UI, Controller"
$sprintBacklog.Range("F10").Value = "Carole"
$sprintBacklog.Range("G10").Value = "Jonas"
$sprintBacklog.Range("I10").Value = 5
$sprintBacklog.Range("L10").Value = "work in progress"

# Row 11 - Stop a Medication
$sprintBacklog.Range("A11").Value = 4.1
$sprintBacklog.Range("B11").Value = 2
$sprintBacklog.Range("C11").Value = "Stop a Medication"
$sprintBacklog.Range("D11").Value = "Stop an existing Medication over the Medicationoverview-List"
$sprintBacklog.Range("E11").Value = "UI, Controller, Modell"
$sprintBacklog.Range("F11").Value = "Michel"
$sprintBacklog.Range("G11").Value = "Marwin"
$sprintBacklog.Range("I11").Value = 10
$sprintBacklog.Range("L11").Value = "work in progress"

# Row 12 - GUI Optimization
$sprintBacklog.Range("A12").Value = 1.4
$sprintBacklog.Range("B12").Value = 2
$sprintBacklog.Range("C12").Value = "GUI Optimization"
$sprintBacklog.Range("D12").Value = "Create a nice SCSS file to beautify our GUI"
$sprintBacklog.Range("E12").Value = "UI"
$sprintBacklog.Range("F12").Value = "Carole"
$sprintBacklog.Range("G12").Value = "Amin"
$sprintBacklog.Range("I12").Value = 10
$sprintBacklog.Range("L12").Value = "work in progress"

$sprintBacklog.Range("F12").Select()
